$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 68
$ws.Range("B68").Value = "Reglage d'un bug de log in sur le server"
$ws.Range("C68").Value = "5/25/2021"
$ws.Range("D68").Value = 20
$ws.Range("E68").Value = "Il était déjà present au moment du livrable"

# Row 69
$ws.Range("B69").Value = "Buton supprimer a continuer de coder"
$ws.Range("C69").Value = "5/25/2021"
$ws.Range("D69").Value = 60

# Row 70
$ws.Range("B70").Value = "Probleme de reperage entre la suppresion d'un event et l'update"
$ws.Range("C70").Value = "5/25/2021"
$ws.Range("D70").Value = 60
$ws.Range("E70").Value = "J'ai du identifier les form individuellement avec des id et faire les redirections adequates."

$ws.Rows.Item(68).RowHeight = 30
$ws.Rows.Item(69).RowHeight = 30
$ws.Rows.Item(70).RowHeight = 30

$ws.Range("E70").Select()
